$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add a new row 5, carrying over the same formatting as row 4
$ws.Range("A4:C4").Copy()
$ws.Range("A5:C5").PasteSpecial(-4122)

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Полет"
$ws.Range("C5").Value = 294

# Move the active selection as recorded in the saved workbook
$ws.Range("F6").Select()
